# Update "想去人数" (interest count, column F) values for a handful of
# events across the "展览" (Exhibition), "演出" (Performance) and
# "全部类型" (All types) sheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

$changes = @{
    "展览"     = @{ "F4" = 1134; "F7" = 241; "F14" = 12853; "F16" = 5278 }
    "演出"     = @{ "F2" = 111 }
    "全部类型" = @{ "F4" = 1134; "F7" = 241; "F14" = 12853; "F15" = 111; "F18" = 5278 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellChanges = $changes[$sheetName]
    foreach ($cellAddr in $cellChanges.Keys) {
        $ws.Range($cellAddr).Value = $cellChanges[$cellAddr]
    }
}
